# Scheduled runner: refresh Universalis market-price snapshots for the Leve
# profitability tables (currentAveragePrice*, LevePrice*, LeveProfit* -> cols H:N).
# Per job-class worksheet, only the rows whose linked marketable item saw a
# price move since the last run are touched.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")

# row 28
$ws.Range("H28").Value = 2662.2307
$ws.Range("I28").Value = 2515.7144
$ws.Range("J28").Value = 2833.1667
$ws.Range("K28").Value = 2515.7144
$ws.Range("L28").Value = 2833.1667
$ws.Range("M28").Value = -2030.7144
$ws.Range("N28").Value = -3803.1667

# row 32
$ws.Range("H32").Value = 2186.2856
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 2050.6667
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 2050.6667
$ws.Range("M32").Value = -2674
$ws.Range("N32").Value = -2702.6667

# row 43
$ws.Range("H43").Value = 14725.125
$ws.Range("I43").Value = 3933
$ws.Range("J43").Value = 21200.4
$ws.Range("K43").Value = 3933
$ws.Range("L43").Value = 21200.4
$ws.Range("M43").Value = -3864
$ws.Range("N43").Value = -21338.4

# row 64
$ws.Range("H64").Value = 3535.7144
$ws.Range("I64").Value = 3750
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 3750
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -3502
$ws.Range("N64").Value = -3996

# row 67
$ws.Range("H67").Value = 3535.7144
$ws.Range("I67").Value = 3750
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 3750
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -2892
$ws.Range("N67").Value = -5216

# row 74
$ws.Range("H74").Value = 4117.5884
$ws.Range("I74").Value = 5249.5
$ws.Range("J74").Value = 3966.6667
$ws.Range("K74").Value = 5249.5
$ws.Range("L74").Value = 3966.6667
$ws.Range("M74").Value = -4313.5
$ws.Range("N74").Value = -5838.6667

# row 77
$ws.Range("H77").Value = 4117.5884
$ws.Range("I77").Value = 5249.5
$ws.Range("J77").Value = 3966.6667
$ws.Range("K77").Value = 26247.5
$ws.Range("L77").Value = 19833.3335
$ws.Range("M77").Value = -21567.5
$ws.Range("N77").Value = -29193.3335

# row 96
$ws.Range("H96").Value = 1051.909
$ws.Range("I96").Value = 1042.6666
$ws.Range("J96").Value = 1063
$ws.Range("K96").Value = 3127.9998
$ws.Range("L96").Value = 3189
$ws.Range("M96").Value = -1754.9998
$ws.Range("N96").Value = -5935

# row 106
$ws.Range("H106").Value = 281382.78
$ws.Range("I106").Value = 419909.16
$ws.Range("J106").Value = 4330
$ws.Range("K106").Value = 419909.16
$ws.Range("L106").Value = 4330
$ws.Range("M106").Value = -419278.16
$ws.Range("N106").Value = -5592

# row 107
$ws.Range("H107").Value = 11184.667
$ws.Range("I107").Value = 13221
$ws.Range("J107").Value = 1003
$ws.Range("K107").Value = 13221
$ws.Range("L107").Value = 1003
$ws.Range("M107").Value = -11301
$ws.Range("N107").Value = -4843

# row 121
$ws.Range("H121").Value = 1420.8334
$ws.Range("I121").Value = 1000
$ws.Range("J121").Value = 1505
$ws.Range("K121").Value = 3000
$ws.Range("L121").Value = 4515
$ws.Range("M121").Value = -1253
$ws.Range("N121").Value = -8009

# row 134
$ws.Range("H134").Value = 49461.133
$ws.Range("I134").Value = 20709
$ws.Range("J134").Value = 56649.168
$ws.Range("K134").Value = 20709
$ws.Range("L134").Value = 56649.168
$ws.Range("M134").Value = -15639
$ws.Range("N134").Value = -66789.16800000001

# row 139
$ws.Range("H139").Value = 35300.625
$ws.Range("I139").Value = 20709
$ws.Range("J139").Value = 44055.6
$ws.Range("K139").Value = 20709
$ws.Range("L139").Value = 44055.6
$ws.Range("M139").Value = -15569
$ws.Range("N139").Value = -54335.6

# row 141
$ws.Range("H141").Value = 1885.4166
$ws.Range("I141").Value = 1312.5
$ws.Range("J141").Value = 4750
$ws.Range("K141").Value = 3937.5
$ws.Range("L141").Value = 14250
$ws.Range("M141").Value = 1242.5
$ws.Range("N141").Value = -24610


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")

# row 63
$ws.Range("H63").Value = 2241
$ws.Range("I63").Value = 2241
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2241
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1555
$ws.Range("N63").ClearContents()

# row 66
$ws.Range("H66").Value = 2241
$ws.Range("I66").Value = 2241
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11205
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7773
$ws.Range("N66").ClearContents()

# row 74
$ws.Range("H74").Value = 17245560
$ws.Range("I74").Value = 26318242
$ws.Range("J74").Value = 7462.8
$ws.Range("K74").Value = 26318242
$ws.Range("L74").Value = 7462.8
$ws.Range("M74").Value = -26317368
$ws.Range("N74").Value = -9210.799999999999

# row 77
$ws.Range("H77").Value = 17245560
$ws.Range("I77").Value = 26318242
$ws.Range("J77").Value = 7462.8
$ws.Range("K77").Value = 131591210
$ws.Range("L77").Value = 37314
$ws.Range("M77").Value = -131586842
$ws.Range("N77").Value = -46050

# row 97
$ws.Range("H97").Value = 8391.9375
$ws.Range("I97").Value = 8260.23
$ws.Range("J97").Value = 8962.666999999999
$ws.Range("K97").Value = 8260.23
$ws.Range("L97").Value = 8962.666999999999
$ws.Range("M97").Value = -7764.23
$ws.Range("N97").Value = -9954.666999999999


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")

# row 94
$ws.Range("H94").Value = 1333.2307
$ws.Range("I94").Value = 1244.4286
$ws.Range("J94").Value = 1436.8334
$ws.Range("K94").Value = 1244.4286
$ws.Range("L94").Value = 1436.8334
$ws.Range("M94").Value = -793.4286
$ws.Range("N94").Value = -2338.8334


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")

# row 58
$ws.Range("H58").Value = 1630.4865
$ws.Range("I58").Value = 684
$ws.Range("J58").Value = 3602.3333
$ws.Range("K58").Value = 684
$ws.Range("L58").Value = 3602.3333
$ws.Range("M58").Value = -481
$ws.Range("N58").Value = -4008.3333

# row 136
$ws.Range("H136").Value = 1630.4865
$ws.Range("I136").Value = 684
$ws.Range("J136").Value = 3602.3333
$ws.Range("K136").Value = 2052
$ws.Range("L136").Value = 10806.9999
$ws.Range("M136").Value = 498
$ws.Range("N136").Value = -15906.9999


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")

# row 97
$ws.Range("H97").Value = 2013.3334
$ws.Range("I97").Value = 1968.3158
$ws.Range("J97").Value = 2184.4
$ws.Range("K97").Value = 1968.3158
$ws.Range("L97").Value = 2184.4
$ws.Range("M97").Value = -1472.3158
$ws.Range("N97").Value = -3176.4


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")

# row 82
$ws.Range("H82").Value = 2335.2727
$ws.Range("I82").Value = 1951.3334
$ws.Range("J82").Value = 2796
$ws.Range("K82").Value = 1951.3334
$ws.Range("L82").Value = 2796
$ws.Range("M82").Value = -1590.3334
$ws.Range("N82").Value = -3518

# row 85
$ws.Range("H85").Value = 2335.2727
$ws.Range("I85").Value = 1951.3334
$ws.Range("J85").Value = 2796
$ws.Range("K85").Value = 1951.3334
$ws.Range("L85").Value = 2796
$ws.Range("M85").Value = -703.3334
$ws.Range("N85").Value = -5292

# row 93
$ws.Range("H93").Value = 967
$ws.Range("I93").Value = 834
$ws.Range("J93").Value = 1299.5
$ws.Range("K93").Value = 834
$ws.Range("L93").Value = 1299.5
$ws.Range("M93").Value = 414
$ws.Range("N93").Value = -3795.5

# row 139
$ws.Range("H139").Value = 59920.57
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 59920.57
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 59920.57
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -70200.57000000001


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")

# row 62
$ws.Range("H62").Value = 13150.3
$ws.Range("I62").Value = 5933.3335
$ws.Range("J62").Value = 16243.286
$ws.Range("K62").Value = 5933.3335
$ws.Range("L62").Value = 16243.286
$ws.Range("M62").Value = -5309.3335
$ws.Range("N62").Value = -17491.286

# row 65
$ws.Range("H65").Value = 13150.3
$ws.Range("I65").Value = 5933.3335
$ws.Range("J65").Value = 16243.286
$ws.Range("K65").Value = 29666.6675
$ws.Range("L65").Value = 81216.42999999999
$ws.Range("M65").Value = -26546.6675
$ws.Range("N65").Value = -87456.42999999999

# row 136
$ws.Range("H136").Value = 3629.1428
$ws.Range("I136").Value = 2734
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 8202
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -5652
$ws.Range("N136").Value = -32100
